$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Cells.Item(12, 8).Value = 327.16666
$ws.Cells.Item(12, 9).Value = 234.33333
$ws.Cells.Item(12, 10).Value = 420
$ws.Cells.Item(12, 11).Value = 234.33333
$ws.Cells.Item(12, 12).Value = 420
$ws.Cells.Item(12, 13).Value = -64.33332999999999
$ws.Cells.Item(12, 14).Value = -760

# Row 76
$ws.Cells.Item(76, 8).Value = 7735.6665
$ws.Cells.Item(76, 10).Value = 7702
$ws.Cells.Item(76, 12).Value = 7702
$ws.Cells.Item(76, 14).Value = -8332

# Row 79
$ws.Cells.Item(79, 8).Value = 7735.6665
$ws.Cells.Item(79, 10).Value = 7702
$ws.Cells.Item(79, 12).Value = 7702
$ws.Cells.Item(79, 14).Value = -9886

# Row 80
$ws.Cells.Item(80, 8).Value = 1141.1428
$ws.Cells.Item(80, 9).Value = 831.3333
$ws.Cells.Item(80, 11).Value = 2493.9999
$ws.Cells.Item(80, 13).Value = -1495.9999

# Row 82
$ws.Cells.Item(82, 8).Value = 1036.875
$ws.Cells.Item(82, 9).Value = 1036.875
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 3110.625
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).Value = -2704.625
$ws.Cells.Item(82, 14).ClearContents()

# Row 83
$ws.Cells.Item(83, 8).Value = 1141.1428
$ws.Cells.Item(83, 9).Value = 831.3333
$ws.Cells.Item(83, 11).Value = 7481.9997
$ws.Cells.Item(83, 13).Value = -2489.9997

# Row 85
$ws.Cells.Item(85, 8).Value = 1036.875
$ws.Cells.Item(85, 9).Value = 1036.875
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 3110.625
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = -1706.625
$ws.Cells.Item(85, 14).ClearContents()

# Row 104
$ws.Cells.Item(104, 8).Value = 212
$ws.Cells.Item(104, 9).Value = 212
$ws.Cells.Item(104, 11).Value = 636
$ws.Cells.Item(104, 13).Value = 1111


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Cells.Item(61, 8).Value = 1853.3636
$ws.Cells.Item(61, 9).Value = 1738.7
$ws.Cells.Item(61, 11).Value = 1738.7
$ws.Cells.Item(61, 13).Value = -1526.7

# Row 88
$ws.Cells.Item(88, 8).Value = 4450.5
$ws.Cells.Item(88, 9).Value = 1700
$ws.Cells.Item(88, 10).Value = 4843.4287
$ws.Cells.Item(88, 11).Value = 1700
$ws.Cells.Item(88, 12).Value = 4843.4287
$ws.Cells.Item(88, 13).Value = -1294
$ws.Cells.Item(88, 14).Value = -5655.4287

# Row 91
$ws.Cells.Item(91, 8).Value = 4450.5
$ws.Cells.Item(91, 9).Value = 1700
$ws.Cells.Item(91, 10).Value = 4843.4287
$ws.Cells.Item(91, 11).Value = 1700
$ws.Cells.Item(91, 12).Value = 4843.4287
$ws.Cells.Item(91, 13).Value = -296
$ws.Cells.Item(91, 14).Value = -7651.4287

# Row 136
$ws.Cells.Item(136, 8).Value = 1853.3636
$ws.Cells.Item(136, 9).Value = 1738.7
$ws.Cells.Item(136, 11).Value = 5216.1
$ws.Cells.Item(136, 13).Value = -2666.1


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 22
$ws.Cells.Item(22, 8).Value = 5828.857
$ws.Cells.Item(22, 10).Value = 13479.667
$ws.Cells.Item(22, 12).Value = 13479.667
$ws.Cells.Item(22, 14).Value = -13825.667

# Row 86
$ws.Cells.Item(86, 8).Value = 4978.4707
$ws.Cells.Item(86, 9).Value = 4386
$ws.Cells.Item(86, 10).Value = 5824.857
$ws.Cells.Item(86, 11).Value = 4386
$ws.Cells.Item(86, 12).Value = 5824.857
$ws.Cells.Item(86, 13).Value = -3263
$ws.Cells.Item(86, 14).Value = -8070.857

# Row 89
$ws.Cells.Item(89, 8).Value = 4978.4707
$ws.Cells.Item(89, 9).Value = 4386
$ws.Cells.Item(89, 10).Value = 5824.857
$ws.Cells.Item(89, 11).Value = 21930
$ws.Cells.Item(89, 12).Value = 29124.285
$ws.Cells.Item(89, 13).Value = -16314
$ws.Cells.Item(89, 14).Value = -40356.285

# Row 97
$ws.Cells.Item(97, 8).Value = 3428
$ws.Cells.Item(97, 9).Value = 3428
$ws.Cells.Item(97, 11).Value = 3428
$ws.Cells.Item(97, 13).Value = -2437

# Row 100
$ws.Cells.Item(100, 8).Value = 17693.555
$ws.Cells.Item(100, 10).Value = 17693.555
$ws.Cells.Item(100, 12).Value = 17693.555
$ws.Cells.Item(100, 14).Value = -19857.555

# Row 135
$ws.Cells.Item(135, 8).Value = 41498.5
$ws.Cells.Item(135, 10).Value = 41498.5
$ws.Cells.Item(135, 12).Value = 41498.5
$ws.Cells.Item(135, 14).Value = -51638.5


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 36
$ws.Cells.Item(36, 8).Value = 1954.2
$ws.Cells.Item(36, 9).Value = 1117.75
$ws.Cells.Item(36, 11).Value = 1117.75
$ws.Cells.Item(36, 13).Value = -729.75

# Row 40
$ws.Cells.Item(40, 8).Value = 1954.2
$ws.Cells.Item(40, 9).Value = 1117.75
$ws.Cells.Item(40, 11).Value = 1117.75
$ws.Cells.Item(40, 13).Value = -957.75

# Row 58
$ws.Cells.Item(58, 8).Value = 3453.4546
$ws.Cells.Item(58, 9).Value = 2598.9375
$ws.Cells.Item(58, 11).Value = 2598.9375
$ws.Cells.Item(58, 13).Value = -2395.9375

# Row 134
$ws.Cells.Item(134, 8).Value = 2936.8462
$ws.Cells.Item(134, 9).Value = 2596
$ws.Cells.Item(134, 11).Value = 7788
$ws.Cells.Item(134, 13).Value = -5253

# Row 136
$ws.Cells.Item(136, 8).Value = 3453.4546
$ws.Cells.Item(136, 9).Value = 2598.9375
$ws.Cells.Item(136, 11).Value = 7796.8125
$ws.Cells.Item(136, 13).Value = -5246.8125


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 5
$ws.Cells.Item(5, 8).Value = 1154.6471
$ws.Cells.Item(5, 9).Value = 785.2727
$ws.Cells.Item(5, 10).Value = 1831.8334
$ws.Cells.Item(5, 11).Value = 2355.8181
$ws.Cells.Item(5, 12).Value = 5495.5002
$ws.Cells.Item(5, 13).Value = -2243.8181
$ws.Cells.Item(5, 14).Value = -5719.5002

# Row 37
$ws.Cells.Item(37, 8).Value = 99919.8
$ws.Cells.Item(37, 10).Value = 99919.8
$ws.Cells.Item(37, 12).Value = 299759.4
$ws.Cells.Item(37, 14).Value = -299983.4

# Row 57
$ws.Cells.Item(57, 8).Value = 2162.6667
$ws.Cells.Item(57, 9).Value = 2162.6667
$ws.Cells.Item(57, 11).Value = 6488.000100000001
$ws.Cells.Item(57, 13).Value = -5929.000100000001

# Row 86
$ws.Cells.Item(86, 8).Value = 1149.25
$ws.Cells.Item(86, 10).Value = 1203.3334
$ws.Cells.Item(86, 12).Value = 3610.0002
$ws.Cells.Item(86, 14).Value = -5982.0002

# Row 89
$ws.Cells.Item(89, 8).Value = 1149.25
$ws.Cells.Item(89, 10).Value = 1203.3334
$ws.Cells.Item(89, 12).Value = 10830.0006
$ws.Cells.Item(89, 14).Value = -22686.0006

# Row 92
$ws.Cells.Item(92, 8).Value = 603
$ws.Cells.Item(92, 10).Value = 638.6
$ws.Cells.Item(92, 12).Value = 1915.8
$ws.Cells.Item(92, 14).Value = -4411.8

# Row 97
$ws.Cells.Item(97, 8).Value = 636
$ws.Cells.Item(97, 9).Value = 240
$ws.Cells.Item(97, 10).Value = 768
$ws.Cells.Item(97, 11).Value = 720
$ws.Cells.Item(97, 12).Value = 2304
$ws.Cells.Item(97, 13).Value = -224
$ws.Cells.Item(97, 14).Value = -3296

# Row 107
$ws.Cells.Item(107, 8).Value = 1003.5
$ws.Cells.Item(107, 10).Value = 658.5
$ws.Cells.Item(107, 12).Value = 1975.5
$ws.Cells.Item(107, 14).Value = -5815.5

# Row 114
$ws.Cells.Item(114, 8).Value = 428
$ws.Cells.Item(114, 9).Value = 428
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 1284
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).Value = 1970
$ws.Cells.Item(114, 14).ClearContents()

# Row 118
$ws.Cells.Item(118, 8).Value = 6666
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 13).ClearContents()

# Row 135
$ws.Cells.Item(135, 8).Value = 1154.6471
$ws.Cells.Item(135, 9).Value = 785.2727
$ws.Cells.Item(135, 10).Value = 1831.8334
$ws.Cells.Item(135, 11).Value = 7067.454299999999
$ws.Cells.Item(135, 12).Value = 16486.5006
$ws.Cells.Item(135, 13).Value = -4532.454299999999
$ws.Cells.Item(135, 14).Value = -21556.5006


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 19
$ws.Cells.Item(19, 8).Value = 4792.3335
$ws.Cells.Item(19, 10).Value = 4786
$ws.Cells.Item(19, 12).Value = 4786
$ws.Cells.Item(19, 14).Value = -5362

# Row 39
$ws.Cells.Item(39, 8).Value = 39999
$ws.Cells.Item(39, 10).Value = 39999
$ws.Cells.Item(39, 12).Value = 39999
$ws.Cells.Item(39, 14).Value = -41063

# Row 132
$ws.Cells.Item(132, 8).Value = 4737.8
$ws.Cells.Item(132, 9).Value = 4672.25
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 14016.75
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -11486.75
$ws.Cells.Item(132, 14).Value = -20060


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Cells.Item(22, 8).Value = 1084
$ws.Cells.Item(22, 10).Value = 1220
$ws.Cells.Item(22, 12).Value = 1220
$ws.Cells.Item(22, 14).Value = -1810

# Row 27
$ws.Cells.Item(27, 8).Value = 1084
$ws.Cells.Item(27, 10).Value = 1220
$ws.Cells.Item(27, 12).Value = 1220
$ws.Cells.Item(27, 14).Value = -1434

# Row 46
$ws.Cells.Item(46, 8).Value = 3556.5217
$ws.Cells.Item(46, 9).Value = 2733.3333
$ws.Cells.Item(46, 10).Value = 4454.5454
$ws.Cells.Item(46, 11).Value = 2733.3333
$ws.Cells.Item(46, 12).Value = 4454.5454
$ws.Cells.Item(46, 13).Value = -2545.3333
$ws.Cells.Item(46, 14).Value = -4830.5454

# Row 82
$ws.Cells.Item(82, 8).Value = 1666.6666

# Row 85
$ws.Cells.Item(85, 8).Value = 1666.6666

# Row 97
$ws.Cells.Item(97, 8).Value = 10466.667
$ws.Cells.Item(97, 10).Value = 10466.667
$ws.Cells.Item(97, 12).Value = 10466.667
$ws.Cells.Item(97, 14).Value = -12448.667

# Row 101
$ws.Cells.Item(101, 8).Value = 12214.167
$ws.Cells.Item(101, 10).Value = 12214.167
$ws.Cells.Item(101, 12).Value = 12214.167
$ws.Cells.Item(101, 14).Value = -18704.167

# Row 121
$ws.Cells.Item(121, 8).Value = 150000
$ws.Cells.Item(121, 10).Value = 150000
$ws.Cells.Item(121, 12).Value = 150000
$ws.Cells.Item(121, 14).Value = -153494

# Row 132
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).ClearContents()


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 100
$ws.Cells.Item(100, 8).Value = 859.6923
$ws.Cells.Item(100, 9).Value = 919.9
$ws.Cells.Item(100, 11).Value = 1839.8
$ws.Cells.Item(100, 13).Value = -1298.8

# Row 122
$ws.Cells.Item(122, 8).Value = 2843.3
$ws.Cells.Item(122, 9).Value = 2934.8572
$ws.Cells.Item(122, 10).Value = 2629.6667
$ws.Cells.Item(122, 11).Value = 8804.571599999999
$ws.Cells.Item(122, 12).Value = 7889.000100000001
$ws.Cells.Item(122, 13).Value = -6354.571599999999
$ws.Cells.Item(122, 14).Value = -12789.0001

# Row 132
$ws.Cells.Item(132, 8).Value = 1284.25
$ws.Cells.Item(132, 9).Value = 1282.7142
$ws.Cells.Item(132, 11).Value = 3848.1426
$ws.Cells.Item(132, 13).Value = -1318.1426

